$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 406.8095
$ws.Range("I33").Value = 447
$ws.Range("J33").Value = 25
$ws.Range("K33").Value = 447
$ws.Range("L33").Value = 25
$ws.Range("M33").Value = -218
$ws.Range("N33").Value = -483

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1648.5264
$ws.Range("I40").Value = 1374
$ws.Range("J40").Value = 1746.5714
$ws.Range("K40").Value = 1374
$ws.Range("L40").Value = 1746.5714
$ws.Range("M40").Value = -1199
$ws.Range("N40").Value = -2096.5714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1547.3158
$ws.Range("I116").Value = 1359.3
$ws.Range("J116").Value = 1756.2222
$ws.Range("K116").Value = 1359.3
$ws.Range("L116").Value = 1756.2222
$ws.Range("M116").Value = 2082.7
$ws.Range("N116").Value = -8640.2222

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1185.6428
$ws.Range("I137").Value = 1050.5
$ws.Range("J137").Value = 1996.5
$ws.Range("K137").Value = 3151.5
$ws.Range("L137").Value = 5989.5
$ws.Range("M137").Value = -601.5
$ws.Range("N137").Value = -11089.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1358.4736
$ws.Range("I2").Value = 1147
$ws.Range("J2").Value = 1816.6666
$ws.Range("K2").Value = 1147
$ws.Range("L2").Value = 1816.6666
$ws.Range("M2").Value = -1034
$ws.Range("N2").Value = -2042.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1358.4736
$ws.Range("I116").Value = 1147
$ws.Range("J116").Value = 1816.6666
$ws.Range("K116").Value = 1147
$ws.Range("L116").Value = 1816.6666
$ws.Range("M116").Value = 1147
$ws.Range("N116").Value = -6404.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 925
$ws.Range("I122").Value = 930.9167
$ws.Range("J122").Value = 889.5
$ws.Range("K122").Value = 2792.7501
$ws.Range("L122").Value = 2668.5
$ws.Range("M122").Value = -342.7501000000002
$ws.Range("N122").Value = -7568.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1358.4736
$ws.Range("I3").Value = 1147
$ws.Range("J3").Value = 1816.6666
$ws.Range("K3").Value = 1147
$ws.Range("L3").Value = 1816.6666
$ws.Range("M3").Value = -1033
$ws.Range("N3").Value = -2044.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 934.4706
$ws.Range("I16").Value = 876.9167
$ws.Range("J16").Value = 1072.6
$ws.Range("K16").Value = 876.9167
$ws.Range("L16").Value = 1072.6
$ws.Range("M16").Value = -589.9167
$ws.Range("N16").Value = -1646.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8775780
$ws.Range("I31").Value = 3061.2
$ws.Range("J31").Value = 18523246
$ws.Range("K31").Value = 3061.2
$ws.Range("L31").Value = 18523246
$ws.Range("M31").Value = -2766.2
$ws.Range("N31").Value = -18523836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8775780
$ws.Range("I34").Value = 3061.2
$ws.Range("J34").Value = 18523246
$ws.Range("K34").Value = 3061.2
$ws.Range("L34").Value = 18523246
$ws.Range("M34").Value = -2859.2
$ws.Range("N34").Value = -18523650

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2080.3125
$ws.Range("I58").Value = 1933.5714
$ws.Range("J58").Value = 2194.4443
$ws.Range("K58").Value = 1933.5714
$ws.Range("L58").Value = 2194.4443
$ws.Range("M58").Value = -1730.5714
$ws.Range("N58").Value = -2600.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1315.5714
$ws.Range("I105").Value = 649.5
$ws.Range("K105").Value = 649.5
$ws.Range("M105").Value = 1097.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 605.4828
$ws.Range("I107").Value = 637.2353000000001
$ws.Range("J107").Value = 560.5
$ws.Range("K107").Value = 637.2353000000001
$ws.Range("L107").Value = 560.5
$ws.Range("M107").Value = 1282.7647
$ws.Range("N107").Value = -4400.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 934.4706
$ws.Range("I113").Value = 876.9167
$ws.Range("J113").Value = 1072.6
$ws.Range("K113").Value = 876.9167
$ws.Range("L113").Value = 1072.6
$ws.Range("M113").Value = 1293.0833
$ws.Range("N113").Value = -5412.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2080.3125
$ws.Range("I136").Value = 1933.5714
$ws.Range("J136").Value = 2194.4443
$ws.Range("K136").Value = 5800.7142
$ws.Range("L136").Value = 6583.3329
$ws.Range("M136").Value = -3250.7142
$ws.Range("N136").Value = -11683.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 511.42856
$ws.Range("I86").Value = 513.3333
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 1539.9999
$ws.Range("L86").Value = 1500
$ws.Range("M86").Value = -353.9999
$ws.Range("N86").Value = -3872

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 511.42856
$ws.Range("I89").Value = 513.3333
$ws.Range("J89").Value = 500
$ws.Range("K89").Value = 4619.9997
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 1308.0003
$ws.Range("N89").Value = -16356

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3018.3333
$ws.Range("I104").Value = 2013
$ws.Range("K104").Value = 6039
$ws.Range("M104").Value = -3418

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18218154
$ws.Range("I70").Value = 21253824
$ws.Range("J70").Value = 4125
$ws.Range("K70").Value = 21253824
$ws.Range("L70").Value = 4125
$ws.Range("M70").Value = -21253554
$ws.Range("N70").Value = -4665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 18218154
$ws.Range("I73").Value = 21253824
$ws.Range("J73").Value = 4125
$ws.Range("K73").Value = 21253824
$ws.Range("L73").Value = 4125
$ws.Range("M73").Value = -21252888
$ws.Range("N73").Value = -5997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 799
$ws.Range("I61").Value = 799
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 799
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -597

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 799
$ws.Range("I113").Value = 799
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 799
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1371

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2680.7
$ws.Range("I132").Value = 3700.75
$ws.Range("J132").Value = 2000.6666
$ws.Range("K132").Value = 11102.25
$ws.Range("L132").Value = 6001.9998
$ws.Range("M132").Value = -8572.25
$ws.Range("N132").Value = -11061.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 335.57144
$ws.Range("I107").Value = 350
$ws.Range("J107").Value = 324.75
$ws.Range("K107").Value = 1050
$ws.Range("L107").Value = 974.25
$ws.Range("M107").Value = 870
$ws.Range("N107").Value = -4814.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 12110.4
$ws.Range("I126").Value = 22500.8
$ws.Range("J126").Value = 1720
$ws.Range("K126").Value = 67502.39999999999
$ws.Range("L126").Value = 5160
$ws.Range("M126").Value = -65032.39999999999
$ws.Range("N126").Value = -10100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4248.028
$ws.Range("I136").Value = 4826.9033
$ws.Range("J136").Value = 659
$ws.Range("K136").Value = 14480.7099
$ws.Range("L136").Value = 1977
$ws.Range("M136").Value = -11930.7099
$ws.Range("N136").Value = -7077
